{"js": "// 1. Update the date/time stamp in the document's Date paragraph.\nconst body = context.document.body;\nconst dateResults = body.search(\"June  21, 2021 (05:45:45 PM)\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"June  21, 2021 (06:21:23 PM)\", \"Replace\");\n}\n\n// 2. Remove the stray \"<! \u2013 TODO: ... \u2013>\" paragraph that precedes the\n//    \"Draw the UML diagram ...\" paragraph, merging the following paragraph\n//    into the (now gone) TODO paragraph's spot while keeping the\n//    \"FirstParagraph\" style that paragraph originally carried.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet todoParagraph = null;\nlet followingParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"TODO: we should use different class here\") !== -1) {\n    todoParagraph = paragraphs.items[i];\n    followingParagraph = paragraphs.items[i + 1];\n    break;\n  }\n}\n\nif (todoParagraph && followingParagraph) {\n  todoParagraph.delete();\n  followingParagraph.style = \"First Paragraph\";\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the date/time stamp in the document's Date paragraph.\n$dateRange = $d.Content\n$find = $dateRange.Find\n$find.ClearFormatting()\n$find.Text = \"June  21, 2021 (05:45:45 PM)\"\n$foundDate = $find.Execute()\nif ($foundDate) {\n    $datePara = $dateRange.Paragraphs.Item(1).Range\n    $datePara.Text = \"June  21, 2021 (06:21:23 PM)\"\n}\n\n# 2. Remove the stray \"<! - TODO: ... ->\" paragraph that precedes the\n#    \"Draw the UML diagram ...\" paragraph. Deleting the paragraph's Range\n#    (including its own end-of-paragraph mark) merges it into the\n#    following paragraph, which then inherits that following paragraph's\n#    own style (\"Body Text\"); re-apply the \"First Paragraph\" style -\n#    originally carried by the removed paragraph - to the merged result.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*TODO: we should use different class here*\") {\n        $p.Range.Delete()\n        $merged = $d.Paragraphs.Item($i)\n        $merged.Style = \"First Paragraph\"\n        break\n    }\n}\n"}
